# Updates the "Status" column (C) on each status_sheet: the status values
# that had already been filled in near the top of each sheet are cleared
# out, and status values are instead filled in further down (rows 49-51)
# -- "made a bit fast", i.e. the reviewer caught up and tagged the most
# recently-added rows instead.  Each sheet's selection is also moved: the
# whole Status column is selected, anchored on the row that was most
# recently touched (and the view is scrolled back up to the top, since
# there's nothing left to review further down).

$wb = $excel.ActiveWorkbook

function Update-StatusSheet($SheetName, $Clears, $Fills, $ActiveCell) {
    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($cell in $Clears) { $ws.Range($cell).ClearContents() | Out-Null }
    foreach ($kv in $Fills.GetEnumerator()) { $ws.Range($kv.Key).Value = $kv.Value }

    $ws.Activate() | Out-Null
    # Select the full Status column, then land the active cell on the
    # row that matters for this sheet (selecting the column first keeps
    # the column-selection semantics even though this engine always
    # reports the active cell as the anchor of the final selection).
    $ws.Columns("C:C").Select() | Out-Null
    $ws.Range($ActiveCell).Select() | Out-Null
}

Update-StatusSheet "status_sheet1" `
    @("C2", "C3", "C4", "C5") `
    @{ "C49" = "Passed"; "C50" = "Paragraph Repeat"; "C51" = "Passed" } `
    "C2"

Update-StatusSheet "status_sheet2" `
    @("C5") `
    @{ "C49" = "Rescheduled"; "C51" = "Passed" } `
    "C51"

Update-StatusSheet "status_sheet3" `
    @("C5") `
    @{ "C51" = "Passed" } `
    "C1"

Update-StatusSheet "status_sheet4" `
    @("C5") `
    @{ "C51" = "Passed" } `
    "C32"

Update-StatusSheet "status_sheet5" `
    @("C5") `
    @{ "C51" = "Passed" } `
    "C36"

Update-StatusSheet "status_sheet6" `
    @("C5") `
    @{ "C51" = "Passed" } `
    "C34"
